$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update the data values on "ACE_landing_page_data" (sheet1) with the June 2024 refreshed figures ---
$ws1.Range("B2").Value = [double]"470.26551030000002"
$ws1.Range("C2").Value = [double]"8923563036"
$ws1.Range("D2").Value = [double]"18975584.73"
$ws1.Range("E2").Value = [double]"0.88557406599999999"
$ws1.Range("F2").Value = [double]"133.29054529999999"
$ws1.Range("G2").Value = [double]"319.75235679999997"
$ws1.Range("H2").Value = [double]"-0.34658558099999998"
$ws1.Range("I2").Value = [double]"3.4281998000000001E-2"
$ws1.Range("J2").Value = [double]"0.58288823599999995"
$ws1.Range("K2").Value = [double]"0.46956072599999998"
$ws1.Range("L2").Value = [double]"-2.3923535999999999E-2"
$ws1.Range("M2").Value = [double]"-0.35154045499999997"
$ws1.Range("N2").Value = [double]"96.950101239999995"
$ws1.Range("O2").Value = [double]"93.290669390000005"
$ws1.Range("B3").Value = [double]"719.70482509999999"
$ws1.Range("C3").Value = [double]"8627785321"
$ws1.Range("D3").Value = [double]"11987949.810000001"
$ws1.Range("E3").Value = [double]"0.60261141299999998"
$ws1.Range("F3").Value = [double]"136.55748320000001"
$ws1.Range("G3").Value = [double]"493.09530489999997"
$ws1.Range("H3").Value = [double]"-0.25295200800000001"
$ws1.Range("I3").Value = [double]"-4.9996459E-2"
$ws1.Range("J3").Value = [double]"0.27167672100000001"
$ws1.Range("K3").Value = [double]"0.25134547499999998"
$ws1.Range("L3").Value = [double]"-8.2098113E-2"
$ws1.Range("M3").Value = [double]"-0.246572025"
$ws1.Range("N3").Value = [double]"93.736622569999994"
$ws1.Range("O3").Value = [double]"58.936990790000003"
$ws1.Range("B4").Value = [double]"963.39837990000001"
$ws1.Range("C4").Value = [double]"9081845434"
$ws1.Range("D4").Value = [double]"9426884.6850000005"
$ws1.Range("E4").Value = [double]"0.48157077700000001"
$ws1.Range("F4").Value = [double]"148.77132850000001"
$ws1.Range("G4").Value = [double]"654.46906869999998"
$ws1.Range("H4").Value = [double]"1.213430842"
$ws1.Range("I4").Value = [double]"-4.1892621999999997E-2"
$ws1.Range("J4").Value = [double]"-0.56713922999999999"
$ws1.Range("K4").Value = [double]"-0.50710518000000004"
$ws1.Range("L4").Value = [double]"8.4551656000000003E-2"
$ws1.Range("M4").Value = [double]"1.219649314"
$ws1.Range("N4").Value = [double]"98.669761249999993"
$ws1.Range("O4").Value = [double]"46.34589106"
$ws1.Range("B5").Value = [double]"435.25117729999999"
$ws1.Range("C5").Value = [double]"9478943216"
$ws1.Range("D5").Value = [double]"21778098.969999999"
$ws1.Range("E5").Value = [double]"0.977025437"
$ws1.Range("F5").Value = [double]"137.1731145"
$ws1.Range("G5").Value = [double]"294.85246369999999"
$ws1.Range("H5").Value = [double]"-2.038008E-3"
$ws1.Range("I5").Value = [double]"1.4498578999999999E-2"
$ws1.Range("J5").Value = [double]"1.6570358E-2"
$ws1.Range("K5").Value = [double]"1.0387739E-2"
$ws1.Range("L5").Value = [double]"6.3813330000000003E-3"
$ws1.Range("M5").Value = [double]"-1.1177139999999999E-3"
$ws1.Range("N5").Value = [double]"102.9840324"
$ws1.Range("O5").Value = [double]"107.0688182"
$ws1.Range("B6").Value = [double]"436.1400342"
$ws1.Range("C6").Value = [double]"9343476090"
$ws1.Range("D6").Value = [double]"21423110.370000001"
$ws1.Range("E6").Value = [double]"0.96698069399999997"
$ws1.Range("F6").Value = [double]"136.30331770000001"
$ws1.Range("G6").Value = [double]"295.18239299999999"
$ws1.Range("H6").Value = [double]"-3.6186771999999999E-2"
$ws1.Range("I6").Value = [double]"1.5122490000000001E-2"
$ws1.Range("J6").Value = [double]"5.3235690000000002E-2"
$ws1.Range("K6").Value = [double]"4.9273702000000003E-2"
$ws1.Range("L6").Value = [double]"-1.7275999999999999E-4"
$ws1.Range("M6").Value = [double]"-3.0874655000000001E-2"
$ws1.Range("N6").Value = [double]"101.512249"
$ws1.Range("O6").Value = [double]"105.32356900000001"
$ws1.Range("B7").Value = [double]"452.51509479999999"
$ws1.Range("C7").Value = [double]"9204284392"
$ws1.Range("D7").Value = [double]"20340281.460000001"
$ws1.Range("E7").Value = [double]"0.92157145699999998"
$ws1.Range("F7").Value = [double]"136.32686949999999"
$ws1.Range("G7").Value = [double]"304.58639269999998"
$ws1.Range("H7").Value = [double]"-3.5869579999999998E-2"
$ws1.Range("I7").Value = [double]"8.9638819999999994E-3"
$ws1.Range("J7").Value = [double]"4.6501449E-2"
$ws1.Range("K7").Value = [double]"4.3322938999999998E-2"
$ws1.Range("L7").Value = [double]"1.2469209E-2"
$ws1.Range("M7").Value = [double]"-3.8898467999999999E-2"

# --- Remove the "UkSATSE" row from the "ANSP" lookup sheet (last row, row 40) ---
$ws2.Rows(40).Delete() | Out-Null

# --- Re-point the active sheet/selection: ANSP becomes the active/visible tab ---
$ws1.Activate()
$ws1.Range("A1:O7").Select() | Out-Null

$ws2.Activate()
$ws2.Range("K15").Select() | Out-Null
